# Applies the "금고Cell 시스템 및 업무소개" etc. 담당업무 labels into column B
# of the IT 유지보수 업무일지 report, rows 7-43, and updates the sheet's
# active selection / scroll position to reflect the end-state of the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Populate column B (담당업무 항목) for rows 7 through 43.
#    The text values are written in row order so that the underlying
#    shared-string table is built up in the same order as the source
#    workbook (string indices 86-92 in first-use order).
# ---------------------------------------------------------------------
$ws.Range("B7").Value = "금고Cell 시스템 및 업무소개"
$ws.Range("B8").Value = "금고Cell 시스템 및 업무소개"
$ws.Range("B9").Value = "보통예금 부서별 계좌번호 관리"
$ws.Range("B10").Value = "2025 보고서 개선"
$ws.Range("B11").Value = "보통예금 부서별 계좌번호 관리"
$ws.Range("B12").Value = "2025 보고서 개선"
$ws.Range("B13").Value = "보통예금 부서별 계좌번호 관리"
$ws.Range("B14").Value = "2025 보고서 개선"
$ws.Range("B15").Value = "보통예금 부서별 계좌번호 관리"
$ws.Range("B16").Value = "2025 보고서 개선"
$ws.Range("B17").Value = "보통예금 부서별 계좌번호 관리"
$ws.Range("B18").Value = "2025 보고서 개선"
$ws.Range("B19").Value = "보통예금 부서별 계좌번호 관리"
$ws.Range("B20").Value = "2025 보고서 개선"
$ws.Range("B21").Value = "보통예금 부서별 계좌번호 관리"
$ws.Range("B22").Value = "2025 보고서 개선"
$ws.Range("B23").Value = "MMDA, 세입세출외현금 거래구분 추가"
$ws.Range("B24").Value = "2025 보고서 개선"
$ws.Range("B25").Value = "MMDA, 세입세출외현금 거래구분 추가"
$ws.Range("B26").Value = "2025 보고서 개선"
$ws.Range("B27").Value = "MMDA, 세입세출외현금 거래구분 추가"
$ws.Range("B28").Value = "2025 보고서 개선"
$ws.Range("B29").Value = "MMDA, 세입세출외현금 거래구분 추가"
$ws.Range("B30").Value = "2025 보고서 개선"
$ws.Range("B31").Value = "웹취약점 처리"
$ws.Range("B32").Value = "2025 보고서 개선"
$ws.Range("B33").Value = "웹취약점 처리"
$ws.Range("B34").Value = "2025 보고서 개선"
$ws.Range("B35").Value = "웹취약점 처리"
$ws.Range("B36").Value = "2025 보고서 개선"
$ws.Range("B37").Value = "웹취약점 처리"
$ws.Range("B38").Value = "2025 보고서 개선"
$ws.Range("B39").Value = "웹취약점 처리"
$ws.Range("B40").Value = "2025 보고서 개선"
$ws.Range("B41").Value = "회계별 잔액 보고서 수정"
$ws.Range("B42").Value = "2025 보고서 개선"
$ws.Range("B43").Value = "금고운용현황 보고서 수정"

# ---------------------------------------------------------------------
# 2. A few of the newly-labelled cells (B13, B15, B17, B19, B21) switch
#    to the "boxed" border style already used by B9/B11 (full thin box)
#    instead of the previous "open bottom" border. Copy that formatting
#    across without touching the values we just set.
# ---------------------------------------------------------------------
$ws.Range("B9").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("B21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Reflect the final cursor position / scroll state left behind by
#    the author after finishing the edits: viewport scrolled down so
#    row 43 is at the top, with B44 (the next blank row) selected.
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 43
$win.ScrollColumn = 1
$ws.Range("B44").Select()
